$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing content but keep cell formatting (header style s="1" stays intact) ---
$ws.Cells.ClearContents()

# --- Re-write header row (row 1) so shared-string indices 0..19 are recreated in original order ---
$headers = @(
    "Sending cluster",
    "Ligand symbol",
    "Receptor symbol",
    "Target cluster",
    "Ligand-expressing cells",
    "Ligand detection rate",
    "Ligand average expression value",
    "Ligand total expression value",
    "Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value",
    "Receptor-expressing cells",
    "Receptor detection rate",
    "Receptor average expression value",
    "Receptor total expression value",
    "Receptor derived specificity of average expression value",
    "Receptor derived specificity of total expression value",
    "Edge average expression weight",
    "Edge total expression weight",
    "Edge average expression derived specificity",
    "Edge total expression derived specificity"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows 2..13 (Sending cluster / Ligand / Receptor / Target cluster + 16 numeric metrics) ---
# Column order below purposely controls the order in which brand-new shared strings
# ("ECs","FAPs","sCs","Efna5","Epha4","M2") are first introduced, so the resulting
# shared string table matches: 20=ECs,21=FAPs,22=sCs,23=Efna5,24=Epha4,25=M2
$colA = @("ECs","ECs","ECs","ECs","FAPs","FAPs","FAPs","FAPs","sCs","sCs","sCs","sCs")
$colB = @("Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5","Efna5")
$colC = @("Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4","Epha4")
$colD = @("ECs","FAPs","M2","sCs","ECs","FAPs","M2","sCs","ECs","FAPs","M2","sCs")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}

# Numeric columns E..T (16 values per row), pipe-delimited, one row of data per line
$numericRows = @(
    "2|0.6666666666666666|0.4223226666666666|1.266968|0.1676547342089268|0.1676547342089269|3|1|7.839997333333334|23.519992|0.3930770090432645|0.3930770090432645|3.311008580472889|29.799077224256|0.06590122147478844|0.06590122147478845",
    "2|0.6666666666666666|0.4223226666666666|1.266968|0.1676547342089268|0.1676547342089269|3|1|11.57455166666667|34.723655|0.5803178185796234|0.5803178185796235|4.888195525337777|43.99375972804|0.09729302963067099|0.09729302963067102",
    "2|0.6666666666666666|0.4223226666666666|1.266968|0.1676547342089268|0.1676547342089269|2|0.6666666666666666|0.04482366666666667|0.134471|0.002247341686300608|0.002247341686300608|0.01893005043644444|0.170370453928|0.0003767774730933699|0.00037677747309337",
    "2|0.6666666666666666|0.4223226666666666|1.266968|0.1676547342089268|0.1676547342089269|3|1|0.4858216666666667|1.457465|0.02435783069081152|0.02435783069081152|0.2051735017911111|1.84656151612|0.004083705630374046|0.004083705630374047",
    "3|1|1.874986333333333|5.624959|0.7443368783435028|0.7443368783435029|3|1|7.839997333333334|23.519992|0.3930770090432645|0.3930770090432645|14.69988785336978|132.298990680328|0.2925817138598643|0.2925817138598644",
    "3|1|1.874986333333333|5.624959|0.7443368783435028|0.7443368783435029|3|1|11.57455166666667|34.723655|0.5803178185796234|0.5803178185796235|21.70212618946056|195.319135705145|0.4319519535286681|0.4319519535286682",
    "3|1|1.874986333333333|5.624959|0.7443368783435028|0.7443368783435029|2|0.6666666666666666|0.04482366666666667|0.134471|0.002247341686300608|0.002247341686300608|0.08404376240988891|0.7563938616890001|0.001672779295352218|0.001672779295352218",
    "3|1|1.874986333333333|5.624959|0.7443368783435028|0.7443368783435029|3|1|0.4858216666666667|1.457465|0.02435783069081152|0.02435783069081152|0.9109089854372222|8.198180868935001|0.01813043165961821|0.01813043165961822",
    "3|1|0.2216933333333333|0.66508|0.08800838744757017|0.08800838744757018|3|1|7.839997333333334|23.519992|0.3930770090432645|0.3930770090432645|1.738075142151111|15.64267627936|0.03459407370861167|0.03459407370861167",
    "3|1|0.2216933333333333|0.66508|0.08800838744757017|0.08800838744757018|3|1|11.57455166666667|34.723655|0.5803178185796234|0.5803178185796235|2.566000940822222|23.0940084674|0.05107283542028423|0.05107283542028425",
    "3|1|0.2216933333333333|0.66508|0.08800838744757017|0.08800838744757018|2|0.6666666666666666|0.04482366666666667|0.134471|0.002247341686300608|0.002247341686300608|0.009937108075555557|0.08943397268|0.0001977849178550196|0.0001977849178550197",
    "3|1|0.2216933333333333|0.66508|0.08800838744757017|0.08800838744757018|3|1|0.4858216666666667|1.457465|0.02435783069081152|0.02435783069081152|0.1077034246888889|0.9693308222|0.002143693400819256|0.002143693400819256"
)

for ($i = 0; $i -lt $numericRows.Length; $i++) {
    $r = $i + 2
    $parts = $numericRows[$i].Split('|')
    for ($j = 0; $j -lt $parts.Length; $j++) {
        $col = 5 + $j   # column E is index 5
        $ws.Cells.Item($r, $col).Value = [double]$parts[$j]
    }
}
